$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1332
$ws1.Range("F4").Value = 1137
$ws1.Range("F13").Value = 74
$ws1.Range("F15").Value = 700
$ws1.Range("F16").Value = 179
$ws1.Range("F21").Value = 160
$ws1.Range("F22").Value = 676
$ws1.Range("F23").Value = 44
$ws1.Range("F24").Value = 651
$ws1.Range("F25").Value = 161
$ws1.Range("F26").Value = 38
$ws1.Range("F29").Value = 162
$ws1.Range("F31").Value = 276

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 25

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1332
$ws4.Range("F5").Value = 1137
$ws4.Range("F15").Value = 74
$ws4.Range("F17").Value = 700
$ws4.Range("F18").Value = 179
$ws4.Range("F26").Value = 25
$ws4.Range("F29").Value = 160
$ws4.Range("F30").Value = 676
$ws4.Range("F31").Value = 44
$ws4.Range("F32").Value = 651
$ws4.Range("F33").Value = 161
$ws4.Range("F34").Value = 38
$ws4.Range("F39").Value = 162
$ws4.Range("F41").Value = 276
